$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 615.8889
$ws.Range("I2").Value = 509
$ws.Range("J2").Value = 990
$ws.Range("K2").Value = 509
$ws.Range("L2").Value = 990
$ws.Range("M2").Value = -396
$ws.Range("N2").Value = -1216
$ws.Range("H40").Value = 2333.3333
$ws.Range("J40").Value = 2333.3333
$ws.Range("L40").Value = 2333.3333
$ws.Range("N40").Value = -2683.3333
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H137").Value = 2465.2295
$ws.Range("I137").Value = 2146.0732
$ws.Range("K137").Value = 6438.219599999999
$ws.Range("M137").Value = -3888.219599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5433.3335
$ws.Range("I63").Value = 1635
$ws.Range("J63").Value = 7332.5
$ws.Range("K63").Value = 1635
$ws.Range("L63").Value = 7332.5
$ws.Range("M63").Value = -949
$ws.Range("N63").Value = -8704.5
$ws.Range("H66").Value = 5433.3335
$ws.Range("I66").Value = 1635
$ws.Range("J66").Value = 7332.5
$ws.Range("K66").Value = 8175
$ws.Range("L66").Value = 36662.5
$ws.Range("M66").Value = -4743
$ws.Range("N66").Value = -43526.5
$ws.Range("H74").Value = 2466.3914
$ws.Range("I74").Value = 1863.375
$ws.Range("K74").Value = 1863.375
$ws.Range("M74").Value = -989.375
$ws.Range("H77").Value = 2466.3914
$ws.Range("I77").Value = 1863.375
$ws.Range("K77").Value = 9316.875
$ws.Range("M77").Value = -4948.875
$ws.Range("H97").Value = 1003.75
$ws.Range("I97").Value = 989.61536
$ws.Range("J97").Value = 1065
$ws.Range("K97").Value = 989.61536
$ws.Range("L97").Value = 1065
$ws.Range("M97").Value = -493.61536
$ws.Range("N97").Value = -2057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1484.5714
$ws.Range("I86").Value = 1484.5714
$ws.Range("K86").Value = 1484.5714
$ws.Range("M86").Value = -361.5714
$ws.Range("H89").Value = 1484.5714
$ws.Range("I89").Value = 1484.5714
$ws.Range("K89").Value = 7422.857
$ws.Range("M89").Value = -1806.857
$ws.Range("H107").Value = 100659.8
$ws.Range("I107").Value = 111510.89
$ws.Range("K107").Value = 111510.89
$ws.Range("M107").Value = -109590.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6799.5747
$ws.Range("I31").Value = 1236.48
$ws.Range("J31").Value = 13121.272
$ws.Range("K31").Value = 1236.48
$ws.Range("L31").Value = 13121.272
$ws.Range("M31").Value = -941.48
$ws.Range("N31").Value = -13711.272
$ws.Range("H34").Value = 6799.5747
$ws.Range("I34").Value = 1236.48
$ws.Range("J34").Value = 13121.272
$ws.Range("K34").Value = 1236.48
$ws.Range("L34").Value = 13121.272
$ws.Range("M34").Value = -1034.48
$ws.Range("N34").Value = -13525.272
$ws.Range("H76").Value = 4000
$ws.Range("I76").Value = 4000
$ws.Range("K76").Value = 4000
$ws.Range("M76").Value = -3685
$ws.Range("H79").Value = 4000
$ws.Range("I79").Value = 4000
$ws.Range("K79").Value = 4000
$ws.Range("M79").Value = -2908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 904.4375
$ws.Range("J5").Value = 1433.5
$ws.Range("L5").Value = 4300.5
$ws.Range("N5").Value = -4524.5
$ws.Range("H31").Value = 1553.4546
$ws.Range("J31").Value = 1408.8
$ws.Range("L31").Value = 4226.4
$ws.Range("N31").Value = -4802.4
$ws.Range("H68").Value = 900.3333
$ws.Range("J68").Value = 775
$ws.Range("L68").Value = 2325
$ws.Range("N68").Value = -3947
$ws.Range("H71").Value = 900.3333
$ws.Range("J71").Value = 775
$ws.Range("L71").Value = 6975
$ws.Range("N71").Value = -15087
$ws.Range("H87").Value = 1642.6666
$ws.Range("I87").Value = 1642.6666
$ws.Range("K87").Value = 4927.9998
$ws.Range("M87").Value = -3679.9998
$ws.Range("H90").Value = 1642.6666
$ws.Range("I90").Value = 1642.6666
$ws.Range("K90").Value = 14783.9994
$ws.Range("M90").Value = -8543.999400000001
$ws.Range("H113").Value = 881.30615
$ws.Range("I113").Value = 667.3103599999999
$ws.Range("J113").Value = 1191.6
$ws.Range("K113").Value = 2001.93108
$ws.Range("L113").Value = 3574.8
$ws.Range("M113").Value = 168.0689200000002
$ws.Range("N113").Value = -7914.799999999999
$ws.Range("H135").Value = 904.4375
$ws.Range("J135").Value = 1433.5
$ws.Range("L135").Value = 12901.5
$ws.Range("N135").Value = -17971.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 36415360
$ws.Range("J80").Value = 59729.285
$ws.Range("L80").Value = 59729.285
$ws.Range("N80").Value = -61725.285
$ws.Range("H83").Value = 36415360
$ws.Range("J83").Value = 59729.285
$ws.Range("L83").Value = 298646.425
$ws.Range("N83").Value = -308630.425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 824.2222
$ws.Range("I55").Value = 647
$ws.Range("J55").Value = 966
$ws.Range("K55").Value = 647
$ws.Range("L55").Value = 966
$ws.Range("M55").Value = -474
$ws.Range("N55").Value = -1312
$ws.Range("H68").Value = 1561.875
$ws.Range("I68").Value = 1507.5
$ws.Range("J68").Value = 1725
$ws.Range("K68").Value = 1507.5
$ws.Range("L68").Value = 1725
$ws.Range("M68").Value = -758.5
$ws.Range("N68").Value = -3223
$ws.Range("H71").Value = 1561.875
$ws.Range("I71").Value = 1507.5
$ws.Range("J71").Value = 1725
$ws.Range("K71").Value = 7537.5
$ws.Range("L71").Value = 8625
$ws.Range("M71").Value = -3793.5
$ws.Range("N71").Value = -16113
$ws.Range("H82").Value = 22728910
$ws.Range("I82").Value = 35716130
$ws.Range("J82").Value = 1275.75
$ws.Range("K82").Value = 35716130
$ws.Range("L82").Value = 1275.75
$ws.Range("M82").Value = -35715769
$ws.Range("N82").Value = -1997.75
$ws.Range("H85").Value = 22728910
$ws.Range("I85").Value = 35716130
$ws.Range("J85").Value = 1275.75
$ws.Range("K85").Value = 35716130
$ws.Range("L85").Value = 1275.75
$ws.Range("M85").Value = -35714882
$ws.Range("N85").Value = -3771.75
$ws.Range("H132").Value = 2905.3333
$ws.Range("I132").Value = 2194.4348
$ws.Range("K132").Value = 6583.3044
$ws.Range("M132").Value = -4053.3044
$ws.Range("H136").Value = 9805800
$ws.Range("I136").Value = 2433.3333
$ws.Range("J136").Value = 15153091
$ws.Range("K136").Value = 7299.999899999999
$ws.Range("L136").Value = 45459273
$ws.Range("M136").Value = -4749.999899999999
$ws.Range("N136").Value = -45464373

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H122").Value = 2396.7585
$ws.Range("I122").Value = 2341.5
$ws.Range("J122").Value = 2662
$ws.Range("K122").Value = 7024.5
$ws.Range("L122").Value = 7986
$ws.Range("M122").Value = -4574.5
$ws.Range("N122").Value = -12886
$ws.Range("H132").Value = 5954563.5
$ws.Range("I132").Value = 2482
$ws.Range("J132").Value = 9805910
$ws.Range("K132").Value = 7446
$ws.Range("L132").Value = 29417730
$ws.Range("M132").Value = -4916
$ws.Range("N132").Value = -29422790
$ws.Range("H136").Value = 3390.8215
$ws.Range("I136").Value = 2757.9
$ws.Range("K136").Value = 8273.700000000001
$ws.Range("M136").Value = -5723.700000000001
